# RPA datasets push 2023-10-14
# Insert a new IPO record at the top of the data table (row 2), shifting all
# existing records down by one row, and drop the oldest record that falls
# off the bottom of the table (former row 17, 넥스틸/하나).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row between the current row 2 and row 3 so the new row
# inherits the plain (unstyled) formatting of the data rows rather than the
# bold header formatting that Rows("2:2").Insert() would copy down.
$ws.Rows("3:3").Insert()

# Shift the former row 2 (2023-09-21 / 두산로보틱스 ...) down into row 3.
$ws.Range("A2:T2").Copy()
$ws.Range("A3").PasteSpecial()

# Write the new, latest IPO entry into row 2. The date-like text values
# (A2, D2, E2) are forced to text so they are stored as shared strings
# like the rest of the date column, instead of being auto-converted to
# Excel date serials.
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "2023-09-25"
$ws.Range("A2").Style = "Normal"

$ws.Range("B2").Value = "에이치엠씨제6호스팩"
$ws.Range("C2").Value = "현대차"

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "2023-10-04"
$ws.Range("D2").Style = "Normal"

$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "2023-10-13"
$ws.Range("E2").Style = "Normal"

$ws.Range("F2").Value = 8000000
$ws.Range("G2").Value = 4000000
$ws.Range("H2").Value = "-"
$ws.Range("I2").Value = 2000
$ws.Range("J2").Value = 2000
$ws.Range("K2").Value = "-"
$ws.Range("L2").Value = 2000
$ws.Range("M2").Value = "-"
$ws.Range("N2").Value = "-"
$ws.Range("O2").Value = 0
$ws.Range("P2").Value = "-"
$ws.Range("Q2").Value = "-"
$ws.Range("R2").Value = "61.15 : 1"
$ws.Range("S2").Value = "-"
$ws.Range("T2").Value = "-"

# Remove the oldest record, which has now been pushed down to row 18.
$ws.Rows("18:18").Delete()
